$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.62005599999999
$ws.Range("H2").Value = 289.860168
$ws.Range("I2").Value = 0.2116037895476247
$ws.Range("J2").Value = 0.2183905833651517
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.006255666666667
$ws.Range("N2").Value = 3.018767
$ws.Range("O2").Value = 0.03991953272530977
$ws.Range("P2").Value = 0.03991953272530977
$ws.Range("Q2").Value = 97.22447886365066
$ws.Range("R2").Value = 875.020309772856
$ws.Range("S2").Value = 0.008447124401645963
$ws.Range("T2").Value = 0.008718050039544666
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.62005599999999
$ws.Range("H3").Value = 289.860168
$ws.Range("I3").Value = 0.2116037895476247
$ws.Range("J3").Value = 0.2183905833651517
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.646551333333332
$ws.Range("N3").Value = 22.939654
$ws.Range("O3").Value = 0.3033491053003703
$ws.Range("P3").Value = 0.3033491053003703
$ws.Range("Q3").Value = 738.8102180335411
$ws.Range("R3").Value = 6649.291962301871
$ws.Range("S3").Value = 0.0641898202374398
$ws.Range("T3").Value = 0.06624858806984471
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.62005599999999
$ws.Range("H4").Value = 289.860168
$ws.Range("I4").Value = 0.2116037895476247
$ws.Range("J4").Value = 0.2183905833651517
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.55429333333333
$ws.Range("N4").Value = 49.66287999999999
$ws.Range("O4").Value = 0.6567313619743199
$ws.Range("P4").Value = 0.6567313619743199
$ws.Range("Q4").Value = 1599.476748907093
$ws.Range("R4").Value = 14395.29074016384
$ws.Range("S4").Value = 0.1389668449085389
$ws.Range("T4").Value = 0.1434239452557624
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 122.3539896666667
$ws.Range("H5").Value = 367.061969
$ws.Range("I5").Value = 0.2679626668787852
$ws.Range("J5").Value = 0.2765570657541026
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.006255666666667
$ws.Range("N5").Value = 3.018767
$ws.Range("O5").Value = 0.03991953272530977
$ws.Range("P5").Value = 0.03991953272530977
$ws.Range("Q5").Value = 123.1193954413581
$ws.Range("R5").Value = 1108.074558972223
$ws.Range("S5").Value = 0.01069694444962895
$ws.Range("T5").Value = 0.01104002883678654
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 122.3539896666667
$ws.Range("H6").Value = 367.061969
$ws.Range("I6").Value = 0.2679626668787852
$ws.Range("J6").Value = 0.2765570657541026
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.646551333333332
$ws.Range("N6").Value = 22.939654
$ws.Range("O6").Value = 0.3033491053003703
$ws.Range("P6").Value = 0.3033491053003703
$ws.Range("Q6").Value = 935.5860628243028
$ws.Range("R6").Value = 8420.274565418724
$ws.Range("S6").Value = 0.08128623525158069
$ws.Range("T6").Value = 0.08389333846100272
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 122.3539896666667
$ws.Range("H7").Value = 367.061969
$ws.Range("I7").Value = 0.2679626668787852
$ws.Range("J7").Value = 0.2765570657541026
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.55429333333333
$ws.Range("N7").Value = 49.66287999999999
$ws.Range("O7").Value = 0.6567313619743199
$ws.Range("P7").Value = 0.6567313619743199
$ws.Range("Q7").Value = 2025.483835445635
$ws.Range("R7").Value = 18229.35451901072
$ws.Range("S7").Value = 0.1759794871775756
$ws.Range("T7").Value = 0.1816236984563133
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 83.74384566666667
$ws.Range("H8").Value = 251.231537
$ws.Range("I8").Value = 0.1834041070557659
$ws.Range("J8").Value = 0.1892864490617203
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.006255666666667
$ws.Range("N8").Value = 3.018767
$ws.Range("O8").Value = 0.03991953272530977
$ws.Range("P8").Value = 0.03991953272530977
$ws.Range("Q8").Value = 84.26771925054211
$ws.Range("R8").Value = 758.409473254879
$ws.Range("S8").Value = 0.007321406253568863
$ws.Range("T8").Value = 0.007556226597777027
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 83.74384566666667
$ws.Range("H9").Value = 251.231537
$ws.Range("I9").Value = 0.1834041070557659
$ws.Range("J9").Value = 0.1892864490617203
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.646551333333332
$ws.Range("N9").Value = 22.939654
$ws.Range("O9").Value = 0.3033491053003703
$ws.Range("P9").Value = 0.3033491053003703
$ws.Range("Q9").Value = 640.3516147409108
$ws.Range("R9").Value = 5763.164532668197
$ws.Range("S9").Value = 0.05563547178377992
$ws.Range("T9").Value = 0.05741987496835699
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 83.74384566666667
$ws.Range("H10").Value = 251.231537
$ws.Range("I10").Value = 0.1834041070557659
$ws.Range("J10").Value = 0.1892864490617203
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.55429333333333
$ws.Range("N10").Value = 49.66287999999999
$ws.Range("O10").Value = 0.6567313619743199
$ws.Range("P10").Value = 0.6567313619743199
$ws.Range("Q10").Value = 1386.320186027395
$ws.Range("R10").Value = 12476.88167424656
$ws.Range("S10").Value = 0.1204472290184171
$ws.Range("T10").Value = 0.1243103474955863
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 111.321218
$ws.Range("H11").Value = 333.963654
$ws.Range("I11").Value = 0.2438002270031519
$ws.Range("J11").Value = 0.2516196610353779
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.006255666666667
$ws.Range("N11").Value = 3.018767
$ws.Range("O11").Value = 0.03991953272530977
$ws.Range("P11").Value = 0.03991953272530977
$ws.Range("Q11").Value = 112.0176064327353
$ws.Range("R11").Value = 1008.158457894618
$ws.Range("S11").Value = 0.009732391140290272
$ws.Range("T11").Value = 0.01004453929303312
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 111.321218
$ws.Range("H12").Value = 333.963654
$ws.Range("I12").Value = 0.2438002270031519
$ws.Range("J12").Value = 0.2516196610353779
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.646551333333332
$ws.Range("N12").Value = 22.939654
$ws.Range("O12").Value = 0.3033491053003703
$ws.Range("P12").Value = 0.3033491053003703
$ws.Range("Q12").Value = 851.2234079261906
$ws.Range("R12").Value = 7661.010671335715
$ws.Range("S12").Value = 0.07395658073343331
$ws.Range("T12").Value = 0.07632859905106433
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 111.321218
$ws.Range("H13").Value = 333.963654
$ws.Range("I13").Value = 0.2438002270031519
$ws.Range("J13").Value = 0.2516196610353779
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 16.55429333333333
$ws.Range("N13").Value = 49.66287999999999
$ws.Range("O13").Value = 0.6567313619743199
$ws.Range("P13").Value = 0.6567313619743199
$ws.Range("Q13").Value = 1842.844096995946
$ws.Range("R13").Value = 16585.59687296352
$ws.Range("S13").Value = 0.1601112551294283
$ws.Range("T13").Value = 0.1652465226912804
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 42.569235
$ws.Range("H14").Value = 85.13847
$ws.Range("I14").Value = 0.09322920951467238
$ws.Range("J14").Value = 0.06414624078364733
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.006255666666667
$ws.Range("N14").Value = 3.018767
$ws.Range("O14").Value = 0.03991953272530977
$ws.Range("P14").Value = 0.03991953272530977
$ws.Range("Q14").Value = 42.835533944415
$ws.Range("R14").Value = 257.01320366649
$ws.Range("S14").Value = 0.003721666480175725
$ws.Range("T14").Value = 0.00256068795816841
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 42.569235
$ws.Range("H15").Value = 85.13847
$ws.Range("I15").Value = 0.09322920951467238
$ws.Range("J15").Value = 0.06414624078364733
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 7.646551333333332
$ws.Range("N15").Value = 22.939654
$ws.Range("O15").Value = 0.3033491053003703
$ws.Range("P15").Value = 0.3033491053003703
$ws.Range("Q15").Value = 325.5078406482299
$ws.Range("R15").Value = 1953.04704388938
$ws.Range("S15").Value = 0.02828099729413664
$ws.Range("T15").Value = 0.01945870475010155
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 42.569235
$ws.Range("H16").Value = 85.13847
$ws.Range("I16").Value = 0.09322920951467238
$ws.Range("J16").Value = 0.06414624078364733
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 16.55429333333333
$ws.Range("N16").Value = 49.66287999999999
$ws.Range("O16").Value = 0.6567313619743199
$ws.Range("P16").Value = 0.6567313619743199
$ws.Range("Q16").Value = 704.7036031655998
$ws.Range("R16").Value = 4228.2216189936
$ws.Range("S16").Value = 0.06122654574036002
$ws.Range("T16").Value = 0.04212684807537738
